$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured min/max values for rows 2 and 3
$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = 9.5

# Match the saved selection state: column K selected (active cell K1)
$ws.Columns("K").Select()
